# follow_map_time.xlsx: change the "type" of two survey fields (date -> text,
# time -> text) on the "survey" sheet, and update the active-cell selection.
# (Commit: "Changed date and time fields for JGI app and verified database
# persistence")

$wb = $excel.ActiveWorkbook

# Best-effort: restore the recorded window geometry from the source commit.
# (Cosmetic only - some hosts do not persist Window geometry.)
try {
    $win = $excel.ActiveWindow
    $win.Left = 3320
    $win.Top = 180
    $win.Width = 25600
    $win.Height = 16060
} catch {
}

$ws = $wb.Worksheets.Item("survey")

# Row 2 = FMT_FOL_date ("Date of the follow"): type date -> text
$ws.Range("C2").Value = "text"

# Row 4 = FMT_time ("time at the start of the 15 minute time interval"): type time -> text
$ws.Range("C4").Value = "text"

# Restore the active-cell selection recorded on the "survey" sheet.
$ws.Activate() | Out-Null
$ws.Range("C5").Select() | Out-Null
